$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 96
$ws.Cells.Item(2, 8).Value = 102
$ws.Cells.Item(3, 6).Value = 41
$ws.Cells.Item(3, 8).Value = 42
$ws.Cells.Item(4, 6).Value = 36
$ws.Cells.Item(4, 8).Value = 48
$ws.Cells.Item(5, 6).Value = 135
$ws.Cells.Item(5, 8).Value = 146
$ws.Cells.Item(6, 6).Value = 45
$ws.Cells.Item(6, 8).Value = 55
$ws.Cells.Item(7, 6).Value = 40
$ws.Cells.Item(7, 8).Value = 44
$ws.Cells.Item(8, 6).Value = 9
$ws.Cells.Item(8, 8).Value = 12
$ws.Cells.Item(9, 6).Value = 10
$ws.Cells.Item(9, 8).Value = 11
$ws.Cells.Item(10, 6).Value = 562
$ws.Cells.Item(10, 8).Value = 657
$ws.Cells.Item(11, 6).Value = 377
$ws.Cells.Item(11, 8).Value = 442
$ws.Cells.Item(12, 6).Value = 611
$ws.Cells.Item(12, 8).Value = 697
$ws.Cells.Item(13, 6).Value = 129
$ws.Cells.Item(13, 8).Value = 163
$ws.Cells.Item(14, 6).Value = 111
$ws.Cells.Item(14, 8).Value = 145
$ws.Cells.Item(15, 6).Value = 129
$ws.Cells.Item(15, 8).Value = 180
$ws.Cells.Item(16, 6).Value = 164
$ws.Cells.Item(16, 8).Value = 212
$ws.Cells.Item(17, 6).Value = 90
$ws.Cells.Item(17, 8).Value = 114
$ws.Cells.Item(18, 6).Value = 45
$ws.Cells.Item(18, 8).Value = 62
$ws.Cells.Item(20, 6).Value = 61
$ws.Cells.Item(20, 8).Value = 98
$ws.Cells.Item(21, 6).Value = 109
$ws.Cells.Item(21, 8).Value = 140
$ws.Cells.Item(22, 6).Value = 146
$ws.Cells.Item(22, 8).Value = 188
$ws.Cells.Item(23, 6).Value = 152
$ws.Cells.Item(23, 8).Value = 204
$ws.Cells.Item(24, 6).Value = 222
$ws.Cells.Item(24, 8).Value = 252
$ws.Cells.Item(25, 6).Value = 253
$ws.Cells.Item(25, 8).Value = 313
$ws.Cells.Item(26, 6).Value = 183
$ws.Cells.Item(26, 8).Value = 208
$ws.Cells.Item(27, 6).Value = 283
$ws.Cells.Item(27, 8).Value = 365
$ws.Cells.Item(28, 6).Value = 155
$ws.Cells.Item(28, 8).Value = 207
$ws.Cells.Item(29, 6).Value = 137
$ws.Cells.Item(29, 8).Value = 178
$ws.Cells.Item(30, 6).Value = 206
$ws.Cells.Item(30, 8).Value = 259
$ws.Cells.Item(31, 6).Value = 46
$ws.Cells.Item(31, 8).Value = 73
$ws.Cells.Item(32, 6).Value = 172
$ws.Cells.Item(32, 8).Value = 210
$ws.Cells.Item(33, 6).Value = 230
$ws.Cells.Item(33, 8).Value = 321
$ws.Cells.Item(34, 6).Value = 230
$ws.Cells.Item(34, 8).Value = 268
$ws.Cells.Item(35, 6).Value = 157
$ws.Cells.Item(35, 8).Value = 184
$ws.Cells.Item(36, 6).Value = 75
$ws.Cells.Item(36, 8).Value = 85
$ws.Cells.Item(37, 6).Value = 149
$ws.Cells.Item(37, 8).Value = 185
$ws.Cells.Item(38, 6).Value = 79
$ws.Cells.Item(38, 8).Value = 96
$ws.Cells.Item(39, 6).Value = 135
$ws.Cells.Item(39, 8).Value = 186
$ws.Cells.Item(40, 6).Value = 218
$ws.Cells.Item(40, 8).Value = 298
$ws.Cells.Item(41, 6).Value = 318
$ws.Cells.Item(41, 8).Value = 410
$ws.Cells.Item(42, 6).Value = 374
$ws.Cells.Item(42, 8).Value = 435
$ws.Cells.Item(43, 6).Value = 107
$ws.Cells.Item(43, 8).Value = 134
$ws.Cells.Item(44, 6).Value = 286
$ws.Cells.Item(44, 8).Value = 354
$ws.Cells.Item(45, 6).Value = 145
$ws.Cells.Item(45, 8).Value = 184
$ws.Cells.Item(46, 6).Value = 300
$ws.Cells.Item(46, 8).Value = 364
$ws.Cells.Item(47, 6).Value = 432
$ws.Cells.Item(47, 8).Value = 524
$ws.Cells.Item(48, 6).Value = 185
$ws.Cells.Item(48, 8).Value = 229
$ws.Cells.Item(49, 6).Value = 218
$ws.Cells.Item(49, 8).Value = 305
$ws.Cells.Item(50, 6).Value = 199
$ws.Cells.Item(50, 8).Value = 272
$ws.Cells.Item(51, 6).Value = 158
$ws.Cells.Item(51, 8).Value = 232
$ws.Cells.Item(52, 6).Value = 25
$ws.Cells.Item(52, 8).Value = 33
